# Updated usage section to include information about the Excel sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("tests")

# Rename / re-case the header row (row 1) to the new, cleaned-up titles.
$ws.Range("A1").Value = "Test Number"
$ws.Range("B1").Value = "Test Folder"
$ws.Range("C1").Value = "Temperature Data"
$ws.Range("D1").Value = "IR/RH Data"
$ws.Range("E1").Value = "GPS Data"
$ws.Range("F1").Value = "Test Date"
$ws.Range("G1").Value = "Testing Route"
$ws.Range("H1").Value = "Cart"
$ws.Range("I1").Value = "Set Up"
$ws.Range("J1").Value = "Temperature Configuration"
$ws.Range("K1").Value = "IR/RH Configuration"

# Update the active selection shown in the sheet view.
$ws.Range("L5").Select()
